$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.882.14'
$ws.Range('E2').Value = '  -0.66%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.099.80'
$ws.Range('E3').Value = '  +2.53%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.79'
$ws.Range('E5').Value = '  -0.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.655'
$ws.Range('E6').Value = '  -1.36%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '54.64'
$ws.Range('E8').Value = '  -3.50%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '59.21'
$ws.Range('E9').Value = '  -0.94%  '
$ws.Range('E10').Value = '  -3.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0767'
$ws.Range('E11').Value = '  -1.87%  '
$ws.Range('E12').Value = '  +1.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.943'
$ws.Range('E13').Value = '  +6.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.13'
$ws.Range('E14').Value = '  -7.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.407.15'
$ws.Range('E15').Value = '  +2.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.51'
$ws.Range('E16').Value = '  -3.61%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.072.59'
$ws.Range('E17').Value = '  +0.98%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.842.46'
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.30'
$ws.Range('E19').Value = '  -5.85%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.83'
$ws.Range('E20').Value = '  -2.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0884'
$ws.Range('E21').Value = '  -1.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.49'
$ws.Range('E22').Value = '  +1.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '239.61'
$ws.Range('E23').Value = '  +1.16%  '
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('E25').Value = '  -3.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.73'
$ws.Range('E26').Value = '  +1.44%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.17'
$ws.Range('E27').Value = '  -0.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '167.25'
$ws.Range('E28').Value = '  -1.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '21.05'
$ws.Range('E29').Value = '  +4.55%  '
$ws.Range('E30').Value = '  -1.70%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.29'
$ws.Range('E31').Value = '  +6.88%  '
$ws.Range('E32').Value = '  +1.31%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.74'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0609'
$ws.Range('E34').Value = '  -1.74%  '
$ws.Range('E35').Value = '  +8.88%  '
$ws.Range('E36').Value = '  +0.23%  '
$ws.Range('E37').Value = '  +3.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0830'
$ws.Range('E38').Value = '  -5.92%  '
$ws.Range('E39').Value = '  -4.33%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.17'
$ws.Range('E40').Value = '  +1.38%  '
$ws.Range('E41').Value = '  -6.71%  '
$ws.Range('E42').Value = '  -0.82%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0964'
$ws.Range('E43').Value = '  -2.63%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '96.93'
$ws.Range('E44').Value = '  +1.09%  '
$ws.Range('E45').Value = '  -8.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.432.70'
$ws.Range('E46').Value = '  +12.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.75'
$ws.Range('E47').Value = '  +14.42%  '
$ws.Range('E48').Value = '  -7.90%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.46'
$ws.Range('E49').Value = '  +1.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.295.02'
$ws.Range('E51').Value = '  +2.98%  '
